# Updated cryptos list with GitHub Actions: refresh Price (D) and Volume(1h) (E)
# columns for the crypto rows on Sheet1. Values are written as literal text
# (leading "'" forces text for numeric-looking prices) and the style is reset
# to "Normal" afterward so no stray quote-prefix/number-format style sticks
# to the cell (matches original cells, which carry no explicit style).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.653.30"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.04%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.561.99"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.52%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.39%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'210.26"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.52%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("E6").Value = "  -1.09%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "  -0.45%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'25.12"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +5.48%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "  -0.72%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "  -0.35%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "  -0.24%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'1.784.67"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.54%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'1.561.84"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.74%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'28.662.60"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.07%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("E15").Value = "  +0.25%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'3.64"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.01%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'61.22"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.14%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'228.14"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.16%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "  -0.57%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "  -1.02%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = "  -0.44%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "  -0.81%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'9.01"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.58%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "  +0.87%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'151.16"
$ws.Range("D25").Style = "Normal"
$ws.Range("E26").Value = "  -1.13%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "  +0.35%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "  -0.41%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = "  -1.84%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = "  -4.21%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = "  -2.62%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = "  -0.13%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'1.391.38"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.90%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'2.98"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.66%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = "  -4.37%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = "  -1.44%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'2.69"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.93%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "  -2.28%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.0162"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.77%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = "  +2.03%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.519"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.21%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.998"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.36%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.770"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.72%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.0459"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.54%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'64.00"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.75%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'5.23"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.01%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'1.697.38"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.57%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "  -5.34%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'85.01"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.37%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'43.25"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +6.99%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "  +0.70%  "
$ws.Range("E51").Style = "Normal"
